$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This document has 21 paragraphs. The edit rearranges several paragraphs'
# contents (paragraph count and paragraph styles stay identical; only the
# run text inside specific paragraphs changes). We therefore rewrite the
# text of each affected paragraph directly, using the exact content that
# should end up there.
# ---------------------------------------------------------------------------

$NL = [char]11   # manual line break (becomes <w:br/>)

# ---------------------------------------------------------------------------
# Phase 1: paragraph 17 ("Avaliação" bullets) internal value shifts.
# Do this FIRST, while the other paragraphs still hold their original text,
# because some of these literal values (e.g. "Aplicação de 2 provas, P1 e
# P2.") also need to land verbatim in a different paragraph later (14) --
# we don't want that later write to be matched/clobbered by these
# whole-document Find/Replace calls, nor vice versa.
#
# Replace literal values in place so the bold "Método:"/"Critério:"/
# "Norma de recuperação:" runs (and the ListBullet paragraph style) are
# left completely untouched.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Aplicação de 2 provas, P1 e P2.", $true, $false, $false, $false, $false, $true, 1, $false, "TMP_METODO_VALUE", 2) | Out-Null
$d.Content.Find.Execute("A média do período será MP = (P1+2P2)/3. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental).", $true, $false, $false, $false, $false, $true, 1, $false, "TMP_CRITERIO_VALUE", 2) | Out-Null
$d.Content.Find.Execute("A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "TMP_NORMA_VALUE", 2) | Out-Null

$d.Content.Find.Execute("TMP_METODO_VALUE", $true, $false, $false, $false, $false, $true, 1, $false, "A média do período será MP = (P1+2P2)/3. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental).", 2) | Out-Null
$d.Content.Find.Execute("TMP_CRITERIO_VALUE", $true, $false, $false, $false, $false, $true, 1, $false, "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação.", 2) | Out-Null

$bibliografiaText = "1)INCROPERA, Frank P. Fundamentos de transferência de calor e de massa. 6. ed. Rio de Janeiro: LTC. 2013." + $NL + "2)KREITH, Frank; BOHN, Mark S. Princípios de transferência de calor. São Paulo: Pioneira. 2014." + $NL + "3) ÖZISIC, M. Necati. Transferência de calor. Rio de Janeiro: Guanabara Koogan. 1990." + $NL + "4) HOLMAN, J. P. Transferência de calor. São Paulo: McGraw-Hill, 1983."
$d.Content.Find.Execute("TMP_NORMA_VALUE", $true, $false, $false, $false, $false, $true, 1, $false, $bibliografiaText, 2) | Out-Null

# ---------------------------------------------------------------------------
# Phase 2: whole-paragraph rewrites (each target paragraph is addressed by
# its fixed index, so these no longer risk colliding with each other).
# ---------------------------------------------------------------------------

# --- Paragraph 6: under "Objetivos" -> short PT program summary ----------
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "1) Introdução;" + $NL + "2) Modos de transferência de calor;" + $NL + "3) Condução;" + $NL + "4) Equação diferencial da condução;" + $NL + "5) Superfícies estendidas (aletas); " + $NL + "6) Coeficiente convectivo (método empírico);" + $NL + "7) Análise transiente;" + $NL + "8) Projeto de trocadores de calor."

# --- Paragraph 7: under "Objetivos" -> short EN program summary (italic) -
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "1)Introduction; 2) Heat transfer modes; 3) Conduction; 4) Differential equation of conduction; 5) Extended surfaces (fins); 6) Convective coefficient (empirical method); 7) Transient analysis; 8) Design of heat exchangers."
$p7.Range.Italic = 1

# --- Paragraph 9: under "Docente(s) Responsável(eis)" -> objectives text -
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "Disciplina do núcleo de base que analisa os fenômenos envolvidos no transporte de calor e estuda a modelagem matemática que os descreve. Esta disciplina apresenta e discute os conceitos que regem o transporte de energia de modo a promover a sua aprendizagem bem como dos métodos de resolução de problemas quando da utilização do calor em processos produtivos industriais (Operações Unitárias)"

# --- Paragraph 11: under "Programa resumido" -> long PT program ----------
$p11 = $d.Paragraphs.Item(11)
$p11.Range.Text = "1) Introdução: conceitos gerais dos fenômenos de transferência de calor e relação com a termodinâmica. Conservação de energia;" + $NL + "2) Modos de transferência de calor: condução convecção e radiação;" + $NL + "3) Condução: analogia com circuito elétrico em paredes simples e compostas nas geometrias: plana, cilíndrica e esférica;" + $NL + "4) Equação diferencial da condução: condução em regime estacionário. Condução de calor em meios compostos. Sistemas com geração de calor. Condução em regime transiente;" + $NL + "5) Superfícies estendidas (aletas): aletas com área de seção transversal uniforme (reta). Desempenho das aletas. Eficiência global da superfície; " + $NL + "6) Coeficiente convectivo (método empírico): convecção natural e forçada, convecção em escoamento externo, convecção em escoamento interno, correlações experimentais para a determinação do coeficiente de convecção;" + $NL + "7) Análise transiente: parâmetros concentrados e ábacos;" + $NL + "8) Projeto de trocadores de calor: método LMDT."

# --- Paragraph 12: under "Programa resumido" EN -> objectives EN (italic) -
$p12 = $d.Paragraphs.Item(12)
$p12.Range.Text = "Basic discipline that analyses the phenomena involved in heat transport, witch studies mathematical modeling that describes them. This course introduces and discusses the concepts governing the transport of energy in order to promote their learning as well as troubleshooting methods when using the heat in industrial production processes (unit operations)."
$p12.Range.Italic = 1

# --- Paragraph 14: under "Programa" -> new short method paragraph --------
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Text = "Aplicação de 2 provas, P1 e P2."

# --- Paragraph 19: under "Bibliografia" -> docente list ------------------
$p19 = $d.Paragraphs.Item(19)
$p19.Range.Text = "6666306 - Daniela Helena Pelegrine Guimarães"

Write-Output "done"
